$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert a brand-new "2022-Q4" sheet right after "总计", as a copy of the
#    "2022-Q3" sheet (same layout/header/style), then overwrite its data row
#    with the new quarter's figures.
# ---------------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item("总计")
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Copy($null, $wsTotal)

$q4 = $wb.Worksheets.Item(2)
$q4.Name = "2022-Q4"

# Fund code / name (A2,B2,C2) stay identical to the template - only the
# numeric-as-text figures and the rank change.
$q4.Range("D2").Value2 = "'3.92"
$q4.Range("D2").Style = "Normal"

$q4.Range("E2").Value2 = "'94.38"
$q4.Range("E2").Style = "Normal"

$q4.Range("F2").Value2 = "'2.10"
$q4.Range("F2").Style = "Normal"

$q4.Range("G2").Value2 = "'0.0823"
$q4.Range("G2").Style = "Normal"

$q4.Range("H2").Value2 = 2

# ---------------------------------------------------------------------------
# 2. Insert a new row at the top of the "总计" summary table for 2022-Q4,
#    pushing all the existing quarters down by one row.
# ---------------------------------------------------------------------------
$wsTotal.Rows.Item(2).Insert()
$wsTotal.Range("B2:D2").Style = "Normal"

$wsTotal.Range("A9").Copy()
$wsTotal.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$wsTotal.Range("A2").Value2 = 0
$wsTotal.Range("B2").Value2 = "2022-Q4"
$wsTotal.Range("C2").Value2 = 1
$wsTotal.Range("D2").Value2 = 0.08
